$d = $word.ActiveDocument

function Insert-DoneMarker($rng) {
    # $rng is collapsed at the insertion point (end of the text that should be
    # immediately followed by "(done)"). Insert the "()" shell first (plain
    # formatting, matching the surrounding text), then insert "done" between
    # the parens and color only that inner text red - this avoids the red
    # formatting leaking onto the parentheses.
    $rng.Collapse(0)
    $rng.InsertAfter("()")
    $mid = $rng.Start + 1
    $rngMid = $d.Range($mid, $mid)
    $rngMid.InsertAfter("done")
    $rngMid.Font.Color = 255
}

# --- 1) "Automatic stopout – loss of the prior day’s low" -> insert "(done)" right after "low", before ".  " ---
$rng = $d.Content
$rng.Find.Execute("Automatic stopout – loss of the prior day’s low", $false, $false, $false, $false, $false, $true, 1, $false, "", 0) | Out-Null
Insert-DoneMarker $rng

# --- 2) after "5EMA" insert "(done)" ---
$rng = $d.Content
$rng.Find.Execute("Alternate exit at 20 ticks below the 5EMA", $false, $false, $false, $false, $false, $true, 1, $false, "", 0) | Out-Null
Insert-DoneMarker $rng

# --- 3) after "the daily open" insert "(done)" (before ", or if the") ---
$rng = $d.Content
$rng.Find.Execute("VWAP drops below the daily open", $false, $false, $false, $false, $false, $true, 1, $false, "", 0) | Out-Null
Insert-DoneMarker $rng

# --- 4) after "before any of those" insert "(done)" (before the final ".") ---
$rng = $d.Content
$rng.Find.Execute("price loses the low of the measured sequence of n-5, n-4, n-3, n-2, n-1 before any of those", $false, $false, $false, $false, $false, $true, 1, $false, "", 0) | Out-Null
Insert-DoneMarker $rng

# --- 5) Append " Current 5:00 PM CT" after "Open and close time: 8:30 am CT, 15:00 pm CT?" ---
$rng = $d.Content
$rng.Find.Execute("Open and close time: 8:30 am CT, 15:00 pm CT?", $false, $false, $false, $false, $false, $true, 1, $false, "", 0) | Out-Null
$rng.Collapse(0)
$rng.InsertAfter(" Current 5:00 PM CT")
